$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text updates (rich-text cells; collapse to plain text, which is
# the best achievable via this COM surface, but the visible/stored string
# content matches the source diff exactly) ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Column width swap between col E (5) and col H (8) ---
$ws.Columns.Item(5).ColumnWidth = 6.168446
$ws.Columns.Item(8).ColumnWidth = 7.433768

# --- Bulk numeric value updates for the weekly crime-stat table (rows 14-30) ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 11
$ws.Range("G14").Value = 11
$ws.Range("I14").Value = 81
$ws.Range("J14").Value = 90
$ws.Range("K14").Value = -10
$ws.Range("L14").Value = -20.588235294117
$ws.Range("M14").Value = -16.494845360824
$ws.Range("N14").Value = -77.808219178082
$ws.Range("C15").Value = 11
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = 22.222222222222
$ws.Range("F15").Value = 48
$ws.Range("G15").Value = 36
$ws.Range("H15").Value = 33.333333333333
$ws.Range("I15").Value = 390
$ws.Range("J15").Value = 307
$ws.Range("K15").Value = 27.035830618892
$ws.Range("L15").Value = 38.790035587188
$ws.Range("M15").Value = 73.333333333333
$ws.Range("N15").Value = -25.572519083969
$ws.Range("C16").Value = 106
$ws.Range("E16").Value = 8.163265306122
$ws.Range("F16").Value = 376
$ws.Range("G16").Value = 418
$ws.Range("H16").Value = -10.047846889952
$ws.Range("I16").Value = 3521
$ws.Range("J16").Value = 3630
$ws.Range("K16").Value = -3.002754820936
$ws.Range("L16").Value = -0.085130533484
$ws.Range("M16").Value = 12.708066581306
$ws.Range("N16").Value = -69.575736628359
$ws.Range("C17").Value = 181
$ws.Range("D17").Value = 180
$ws.Range("E17").Value = 0.555555555555
$ws.Range("F17").Value = 698
$ws.Range("G17").Value = 710
$ws.Range("H17").Value = -1.69014084507
$ws.Range("I17").Value = 6413
$ws.Range("J17").Value = 6169
$ws.Range("K17").Value = 3.955260171826
$ws.Range("L17").Value = 7.133311059137
$ws.Range("M17").Value = 96.176200672988
$ws.Range("N17").Value = -4.497393894266
$ws.Range("C18").Value = 42
$ws.Range("D18").Value = 56
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 184
$ws.Range("G18").Value = 234
$ws.Range("H18").Value = -21.367521367521
$ws.Range("I18").Value = 2028
$ws.Range("J18").Value = 2131
$ws.Range("K18").Value = -4.833411543876
$ws.Range("L18").Value = -6.629834254143
$ws.Range("M18").Value = -13.885350318471
$ws.Range("N18").Value = -85.11341114292
$ws.Range("C19").Value = 170
$ws.Range("D19").Value = 207
$ws.Range("E19").Value = -17.874396135265
$ws.Range("F19").Value = 734
$ws.Range("G19").Value = 794
$ws.Range("H19").Value = -7.556675062972
$ws.Range("I19").Value = 6654
$ws.Range("J19").Value = 6813
$ws.Range("K19").Value = -2.333773667987
$ws.Range("L19").Value = 14.270994332818
$ws.Range("M19").Value = 97.330960854092
$ws.Range("N19").Value = 21.113942482708
$ws.Range("C20").Value = 72
$ws.Range("D20").Value = 96
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 297
$ws.Range("G20").Value = 366
$ws.Range("H20").Value = -18.852459016393
$ws.Range("I20").Value = 3251
$ws.Range("J20").Value = 3096
$ws.Range("K20").Value = 5.00645994832
$ws.Range("L20").Value = -16.29763130793
$ws.Range("M20").Value = 116.444740346205
$ws.Range("N20").Value = -70.632339656729
$ws.Range("C21").Value = 584
$ws.Range("D21").Value = 648
$ws.Range("E21").Value = -9.876543209876
$ws.Range("F21").Value = 2348
$ws.Range("G21").Value = 2569
$ws.Range("H21").Value = -8.602569093032
$ws.Range("I21").Value = 22338
$ws.Range("J21").Value = 22236
$ws.Range("K21").Value = 0.45871559633
$ws.Range("L21").Value = 2.599669300018
$ws.Range("M21").Value = 60.197934595525
$ws.Range("N21").Value = -54.748399643464
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = -77.777777777777
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 31
$ws.Range("H22").Value = -32.258064516129
$ws.Range("I22").Value = 219
$ws.Range("J22").Value = 249
$ws.Range("K22").Value = -12.048192771084
$ws.Range("L22").Value = 2.33644859813
$ws.Range("M22").Value = -2.232142857142
$ws.Range("C23").Value = 30
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = -6.25
$ws.Range("F23").Value = 121
$ws.Range("G23").Value = 142
$ws.Range("H23").Value = -14.788732394366
$ws.Range("I23").Value = 1134
$ws.Range("J23").Value = 1278
$ws.Range("K23").Value = -11.267605633802
$ws.Range("L23").Value = -13.632901751713
$ws.Range("M23").Value = 46.511627906976
$ws.Range("C24").Value = 342
$ws.Range("D24").Value = 311
$ws.Range("E24").Value = 9.967845659163
$ws.Range("F24").Value = 1341
$ws.Range("G24").Value = 1275
$ws.Range("H24").Value = 5.176470588235
$ws.Range("I24").Value = 12980
$ws.Range("J24").Value = 11751
$ws.Range("K24").Value = 10.458684367287
$ws.Range("L24").Value = -0.658196846777
$ws.Range("M24").Value = 40.735118724926
$ws.Range("C25").Value = 94
$ws.Range("D25").Value = 115
$ws.Range("E25").Value = -18.260869565217
$ws.Range("F25").Value = 414
$ws.Range("G25").Value = 475
$ws.Range("H25").Value = -12.842105263157
$ws.Range("I25").Value = 4267
$ws.Range("J25").Value = 4668
$ws.Range("K25").Value = -8.590402742073
$ws.Range("L25").Value = -23.338124326266
$ws.Range("C26").Value = 234
$ws.Range("D26").Value = 229
$ws.Range("E26").Value = 2.183406113537
$ws.Range("F26").Value = 886
$ws.Range("G26").Value = 875
$ws.Range("H26").Value = 1.257142857142
$ws.Range("I26").Value = 8146
$ws.Range("J26").Value = 8108
$ws.Range("K26").Value = 0.468672915638
$ws.Range("L26").Value = 6.888859729694
$ws.Range("M26").Value = 0.196801968019
$ws.Range("C27").Value = 16
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = 77.777777777777
$ws.Range("F27").Value = 69
$ws.Range("G27").Value = 44
$ws.Range("H27").Value = 56.818181818181
$ws.Range("I27").Value = 499
$ws.Range("J27").Value = 466
$ws.Range("K27").Value = 7.081545064377
$ws.Range("L27").Value = 4.175365344467
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = 37.5
$ws.Range("F28").Value = 92
$ws.Range("G28").Value = 82
$ws.Range("H28").Value = 12.195121951219
$ws.Range("I28").Value = 789
$ws.Range("J28").Value = 865
$ws.Range("K28").Value = -8.78612716763
$ws.Range("L28").Value = 2.60078023407
$ws.Range("C29").Value = 5
$ws.Range("D29").Value = 8
$ws.Range("E29").Value = -37.5
$ws.Range("F29").Value = 32
$ws.Range("G29").Value = 41
$ws.Range("H29").Value = -21.951219512195
$ws.Range("I29").Value = 247
$ws.Range("J29").Value = 322
$ws.Range("K29").Value = -23.291925465838
$ws.Range("L29").Value = -17.666666666666
$ws.Range("M29").Value = -31.578947368421
$ws.Range("N29").Value = -76.158301158301
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = -16.666666666666
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = -33.333333333333
$ws.Range("I30").Value = 200
$ws.Range("J30").Value = 253
$ws.Range("K30").Value = -20.94861660079
$ws.Range("L30").Value = -18.699186991869
$ws.Range("M30").Value = -33.554817275747
$ws.Range("N30").Value = -78.632478632478

# --- Row 31 (Hate Crimes): F31 switches from a numeric cell to a text cell
# showing "0" (reusing the same shared string / style already used by the
# sibling text cells in this row, e.g. C31/D31/E31). G31/H31/L31 stay numeric. ---
$ws.Range("F31").Value = "'0"
$ws.Range("C31").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
$ws.Range("L31").Value = -6.25

# --- Row 33 (Traffic Fatalities): D33 switches to text "0", E33 switches to
# text "***.*" (both reuse styles already present on sibling text cells). ---
$ws.Range("D33").Value = "'0"
$ws.Range("C33").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E33").Value = "***.*"
$ws.Range("E31").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("L33").Value = -38.888888888888
